$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "dd"
$ws.Range("A12").Value = "dvv"
$ws.Range("A13").Value = "wa"

$ws.Range("A13").Select()
